# Aerodynamic features testing in progress
# Re-order the rows of the per-component "WEIGHT ESTIMATION METHODS COMPARISON"
# tables so that each weight-estimation method (and its own Value / Percent Error
# data) appears in its new row position.

$wb = $excel.ActiveWorkbook

# --- FUSELAGE ---------------------------------------------------------
$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("A8").Value = "SADRAEY"
$ws.Range("C8").Value = 6396.0
$ws.Range("D8").Value = 0.24176597168341904
$ws.Range("A9").Value = "TORENBEEK_1976"
$ws.Range("C9").Value = 10802.0
$ws.Range("D9").Value = 69.29511507600442
$ws.Range("A10").Value = "RAYMER"
$ws.Range("C10").Value = 6463.0
$ws.Range("D10").Value = 1.291828248122254
$ws.Range("A12").Value = "NICOLAI_1984"
$ws.Range("C12").Value = 10243.0
$ws.Range("D12").Value = 60.534147724820706
$ws.Range("A13").Value = "TORENBEEK_2013"
$ws.Range("C13").Value = 7891.0
$ws.Range("D13").Value = 23.67226005043056
$ws.Range("A14").Value = "KROO"
$ws.Range("C14").Value = 7149.0
$ws.Range("D14").Value = 12.043212153152714

# --- WING ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("WING")
$ws.Range("A8").Value = "RAYMER"
$ws.Range("C8").Value = 8394.0
$ws.Range("D8").Value = 24.1090220843084
$ws.Range("A9").Value = "TORENBEEK_1982"
$ws.Range("C9").Value = 6631.0
$ws.Range("D9").Value = -1.9577167689958312
$ws.Range("A10").Value = "TORENBEEK_2013"
$ws.Range("C10").Value = 6138.0
$ws.Range("D10").Value = -9.246940963368484
$ws.Range("A11").Value = "KROO"
$ws.Range("C11").Value = 7561.0
$ws.Range("D11").Value = 11.792746721402883

# --- HORIZONTAL TAIL ------------------------------------------------------
$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("A8").Value = "TORENBEEK_1976"
$ws.Range("C8").Value = 52.0
$ws.Range("D8").Value = -92.91327211228821
$ws.Range("A10").Value = "RAYMER"
$ws.Range("C10").Value = 507.0
$ws.Range("D10").Value = -30.904403094809947
$ws.Range("A11").Value = "JENKINSON"
$ws.Range("C11").Value = 700.0
$ws.Range("D11").Value = -4.601739973110383
$ws.Range("A12").Value = "NICOLAI_2013"
$ws.Range("C12").Value = 399.0
$ws.Range("D12").Value = -45.62299178467292
$ws.Range("A14").Value = "KROO"
$ws.Range("C14").Value = 737.0
$ws.Range("D14").Value = 0.4407394854537826
$ws.Range("A15").Value = "ROSKAM"
$ws.Range("C15").Value = 1523.0
$ws.Range("D15").Value = 107.55935717278983

# --- VERTICAL TAIL --------------------------------------------------------
$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("A8").Value = "TORENBEEK_1976"
$ws.Range("C8").Value = 124.0
$ws.Range("D8").Value = -83.10087965237955
$ws.Range("A10").Value = "RAYMER"
$ws.Range("C10").Value = 180.0
$ws.Range("D10").Value = -75.46901885022838
$ws.Range("A11").Value = "JENKINSON"
$ws.Range("C11").Value = 502.0
$ws.Range("D11").Value = -31.58581923785916
$ws.Range("A12").Value = "HOWE"
$ws.Range("C12").Value = 1145.0
$ws.Range("D12").Value = 56.04429675826945
$ws.Range("A13").Value = "KROO"
$ws.Range("C13").Value = 488.0
$ws.Range("D13").Value = -33.49378443839695

# --- POWER PLANT (ENGINE 1 / ENGINE 2) -------------------------------------
$ws = $wb.Worksheets.Item("POWER PLANT")
$ws.Range("A11").Value = "TORENBEEK_1976"
$ws.Range("C11").Value = 2954.0
$ws.Range("D11").Value = 11.558495400119368
$ws.Range("A12").Value = "KUNDU"
$ws.Range("C12").Value = 3265.0
$ws.Range("D12").Value = 23.303482559712165
$ws.Range("A18").Value = "TORENBEEK_1976"
$ws.Range("C18").Value = 2954.0
$ws.Range("D18").Value = 11.558495400119368
$ws.Range("A19").Value = "KUNDU"
$ws.Range("C19").Value = 3265.0
$ws.Range("D19").Value = 23.303482559712165
